$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Paragraph 4 ("...trova le prenotazioni con importo superiore a 150 e
# caparra inferiore a 50") was split across 3 runs ("...trova ", "le",
# " prenotazioni..."). Re-run a Find & Replace over the (unchanged) full
# text so the engine normalizes/merges it back into a single run.
$d.Content.Find.Execute(
    "Con riferimento al database IFTS scrivi la query che trova le prenotazioni con importo superiore a 150 e caparra inferiore a 50",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Con riferimento al database IFTS scrivi la query che trova le prenotazioni con importo superiore a 150 e caparra inferiore a 50",
    2) | Out-Null

# --- Change 2 -------------------------------------------------------------
# Paragraph 5 (the "regione di provenienza" question) gets struck through -
# both the paragraph mark and its run carry <w:strike/>.
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Font.StrikeThrough = 1

# --- Change 3 -------------------------------------------------------------
# Paragraph 7 ("...per gli hotel 2 stelle...") gets re-split so the digit
# "2" sits alone in its own run (as it does in the authored doc, where it
# is additionally bracketed by grammar-check <w:proofErr> marks).
$p7 = $d.Paragraphs.Item(7)

# normalize back to a single run first (same trick as change 1)
$d.Content.Find.Execute(
    "Con riferimento al database IFTS scrivi la query che trova le date di arrivo per gli hotel 2 stelle ordinate in ordine cronologico dalla più recente alla più vecchia",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Con riferimento al database IFTS scrivi la query che trova le date di arrivo per gli hotel 2 stelle ordinate in ordine cronologico dalla più recente alla più vecchia",
    2) | Out-Null

# isolate the "2" into its own run by toggling a character property on
# just that sub-range (and immediately reverting it) which forces the run
# to split at its boundaries without changing the visible formatting.
$numRange = $d.Range($p7.Range.Start, $p7.Range.End)
$numRange.Find.Execute("2", $false) | Out-Null
$numRange.Bold = 1
$numRange.Bold = 0
